$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: update existing A3, and add new A4/A5/A6 entries.
# Order matters for shared-string append order: RO.ACT.004 must be
# interned before RO.ACT.003 so they land at sharedStrings indices 22/23.
$ws.Range("A3").Value = "RO.ACT.001"
$ws.Range("A5").Value = "RO.ACT.004"
$ws.Range("A4").Value = "RO.ACT.003"
$ws.Range("A6").Value = "AD.SEC.002.FON.01"

# Column D: mirror the same four values into rows 8-11.
$ws.Range("D8").Value = "RO.ACT.001"
$ws.Range("D10").Value = "RO.ACT.004"
$ws.Range("D9").Value = "RO.ACT.003"
$ws.Range("D11").Value = "AD.SEC.002.FON.01"

# Update the active selection shown in the sheet view.
$ws.Range("B17").Select()
